# Apply row permutation of columns D,J,K,L,M,N,O,P,Q per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2 <- original row 44
$ws.Range("D2").Value = 44315
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14538
$ws.Range("N2").Value = '$/malla 15 kilos'
$ws.Range("O2").Value = 'Provincia de Quillota'
$ws.Range("P2").Value = 969
$ws.Range("Q2").Value = 15
# row 3 <- original row 73
$ws.Range("D3").Value = 44791
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("N3").Value = '$/malla 15 kilos'
$ws.Range("O3").Value = 'Provincia de Quillota'
$ws.Range("P3").Value = 1000
$ws.Range("Q3").Value = 15
# row 4 <- original row 64
$ws.Range("D4").Value = 44722
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 15500
$ws.Range("M4").Value = 15263
$ws.Range("N4").Value = '$/malla 15 kilos'
$ws.Range("O4").Value = 'Provincia de Quillota'
$ws.Range("P4").Value = 1018
$ws.Range("Q4").Value = 15
# row 5 <- original row 81
$ws.Range("D5").Value = 44721
$ws.Range("J5").Value = 130
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("N5").Value = '$/malla 15 kilos'
$ws.Range("O5").Value = 'Provincia de Quillota'
$ws.Range("P5").Value = 967
$ws.Range("Q5").Value = 15
# row 6 <- original row 34
$ws.Range("D6").Value = 44343
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("N6").Value = '$/malla 15 kilos'
$ws.Range("O6").Value = 'Provincia de Quillota'
$ws.Range("P6").Value = 1000
$ws.Range("Q6").Value = 15
# row 7 <- original row 67
$ws.Range("D7").Value = 44784
$ws.Range("J7").Value = 105
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14476
$ws.Range("N7").Value = '$/malla 15 kilos'
$ws.Range("O7").Value = 'Provincia de Quillota'
$ws.Range("P7").Value = 965
$ws.Range("Q7").Value = 15
# row 8 <- original row 7
$ws.Range("D8").Value = 44438
$ws.Range("J8").Value = 75
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19467
$ws.Range("N8").Value = '$/malla 15 kilos'
$ws.Range("O8").Value = 'Provincia de Quillota'
$ws.Range("P8").Value = 1298
$ws.Range("Q8").Value = 15
# row 9 <- original row 20
$ws.Range("D9").Value = 44309
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 15000
$ws.Range("N9").Value = '$/malla 15 kilos'
$ws.Range("O9").Value = 'Provincia de Quillota'
$ws.Range("P9").Value = 1000
$ws.Range("Q9").Value = 15
# row 10 <- original row 50
$ws.Range("D10").Value = 45043
$ws.Range("J10").Value = 70
$ws.Range("K10").Value = 11500
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11750
$ws.Range("N10").Value = '$/malla 10 kilos'
$ws.Range("O10").Value = 'Provincia de Quillota'
$ws.Range("P10").Value = 1175
$ws.Range("Q10").Value = 10
# row 11 <- original row 82
$ws.Range("D11").Value = 44802
$ws.Range("J11").Value = 73
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15479
$ws.Range("N11").Value = '$/malla 15 kilos'
$ws.Range("O11").Value = 'Provincia de Quillota'
$ws.Range("P11").Value = 1032
$ws.Range("Q11").Value = 15
# row 12 <- original row 95
$ws.Range("D12").Value = 44782
$ws.Range("J12").Value = 55
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = '$/malla 15 kilos'
$ws.Range("O12").Value = 'Provincia de Quillota'
$ws.Range("P12").Value = 1000
$ws.Range("Q12").Value = 15
# row 13 <- original row 46
$ws.Range("D13").Value = 45040
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("N13").Value = '$/malla 10 kilos'
$ws.Range("O13").Value = 'Provincia de Quillota'
$ws.Range("P13").Value = 1200
$ws.Range("Q13").Value = 10
# row 14 <- original row 79
$ws.Range("D14").Value = 44804
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 15000
$ws.Range("N14").Value = '$/malla 15 kilos'
$ws.Range("O14").Value = 'Provincia de Quillota'
$ws.Range("P14").Value = 1000
$ws.Range("Q14").Value = 15
# row 15 <- original row 62
$ws.Range("D15").Value = 45015
$ws.Range("J15").Value = 73
$ws.Range("K15").Value = 12500
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12740
$ws.Range("N15").Value = '$/malla 10 kilos'
$ws.Range("O15").Value = 'Provincia de Quillota'
$ws.Range("P15").Value = 1274
$ws.Range("Q15").Value = 10
# row 16 <- original row 83
$ws.Range("D16").Value = 44771
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("N16").Value = '$/malla 15 kilos'
$ws.Range("O16").Value = 'Provincia de Quillota'
$ws.Range("P16").Value = 1000
$ws.Range("Q16").Value = 15
# row 17 <- original row 11
$ws.Range("D17").Value = 44322
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("N17").Value = '$/malla 15 kilos'
$ws.Range("O17").Value = 'Provincia de Quillota'
$ws.Range("P17").Value = 967
$ws.Range("Q17").Value = 15
# row 18 <- original row 90
$ws.Range("D18").Value = 45041
$ws.Range("J18").Value = 65
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("N18").Value = '$/malla 10 kilos'
$ws.Range("O18").Value = 'Provincia de Quillota'
$ws.Range("P18").Value = 1200
$ws.Range("Q18").Value = 10
# row 19 <- original row 6
$ws.Range("D19").Value = 44319
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 15000
$ws.Range("N19").Value = '$/malla 15 kilos'
$ws.Range("O19").Value = 'Provincia de Quillota'
$ws.Range("P19").Value = 1000
$ws.Range("Q19").Value = 15
# row 20 <- original row 92
$ws.Range("D20").Value = 44838
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("N20").Value = '$/malla 15 kilos'
$ws.Range("O20").Value = 'Provincia de Quillota'
$ws.Range("P20").Value = 1033
$ws.Range("Q20").Value = 15
# row 21 <- original row 91
$ws.Range("D21").Value = 44748
$ws.Range("J21").Value = 73
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 16000
$ws.Range("M21").Value = 15521
$ws.Range("N21").Value = '$/malla 15 kilos'
$ws.Range("O21").Value = 'Provincia de Quillota'
$ws.Range("P21").Value = 1035
$ws.Range("Q21").Value = 15
# row 22 <- original row 53
$ws.Range("D22").Value = 44719
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 15000
$ws.Range("N22").Value = '$/malla 15 kilos'
$ws.Range("O22").Value = 'Provincia de Quillota'
$ws.Range("P22").Value = 1000
$ws.Range("Q22").Value = 15
# row 24 <- original row 66
$ws.Range("D24").Value = 44727
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("N24").Value = '$/malla 15 kilos'
$ws.Range("O24").Value = 'Provincia de Quillota'
$ws.Range("P24").Value = 1000
$ws.Range("Q24").Value = 15
# row 25 <- original row 77
$ws.Range("D25").Value = 45036
$ws.Range("J25").Value = 105
$ws.Range("K25").Value = 11500
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = 11762
$ws.Range("N25").Value = '$/malla 10 kilos'
$ws.Range("O25").Value = 'Calera'
$ws.Range("P25").Value = 1176
$ws.Range("Q25").Value = 10
# row 26 <- original row 65
$ws.Range("D26").Value = 44811
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 16000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 16000
$ws.Range("N26").Value = '$/malla 15 kilos'
$ws.Range("O26").Value = 'Provincia de Quillota'
$ws.Range("P26").Value = 1067
$ws.Range("Q26").Value = 15
# row 27 <- original row 101
$ws.Range("D27").Value = 44720
$ws.Range("J27").Value = 85
$ws.Range("K27").Value = 15000
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 15529
$ws.Range("N27").Value = '$/malla 15 kilos'
$ws.Range("O27").Value = 'Provincia de Quillota'
$ws.Range("P27").Value = 1035
$ws.Range("Q27").Value = 15
# row 28 <- original row 104
$ws.Range("D28").Value = 44847
$ws.Range("J28").Value = 105
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 16000
$ws.Range("M28").Value = 15524
$ws.Range("N28").Value = '$/malla 15 kilos'
$ws.Range("O28").Value = 'Provincia de Quillota'
$ws.Range("P28").Value = 1035
$ws.Range("Q28").Value = 15
# row 29 <- original row 24
$ws.Range("D29").Value = 44753
$ws.Range("J29").Value = 80
$ws.Range("K29").Value = 15000
$ws.Range("L29").Value = 16000
$ws.Range("M29").Value = 15500
$ws.Range("N29").Value = '$/malla 15 kilos'
$ws.Range("O29").Value = 'Provincia de Quillota'
$ws.Range("P29").Value = 1033
$ws.Range("Q29").Value = 15
# row 30 <- original row 38
$ws.Range("D30").Value = 44795
$ws.Range("J30").Value = 56
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = 15000
$ws.Range("N30").Value = '$/malla 15 kilos'
$ws.Range("O30").Value = 'Provincia de Quillota'
$ws.Range("P30").Value = 1000
$ws.Range("Q30").Value = 15
# row 31 <- original row 71
$ws.Range("D31").Value = 44329
$ws.Range("J31").Value = 35
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 15000
$ws.Range("N31").Value = '$/malla 15 kilos'
$ws.Range("O31").Value = 'Provincia de Quillota'
$ws.Range("P31").Value = 1000
$ws.Range("Q31").Value = 15
# row 32 <- original row 69
$ws.Range("D32").Value = 44340
$ws.Range("J32").Value = 47
$ws.Range("K32").Value = 14000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 14000
$ws.Range("N32").Value = '$/malla 15 kilos'
$ws.Range("O32").Value = 'Provincia de Quillota'
$ws.Range("P32").Value = 933
$ws.Range("Q32").Value = 15
# row 33 <- original row 25
$ws.Range("D33").Value = 44746
$ws.Range("J33").Value = 103
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 16000
$ws.Range("M33").Value = 15563
$ws.Range("N33").Value = '$/malla 15 kilos'
$ws.Range("O33").Value = 'Provincia de Quillota'
$ws.Range("P33").Value = 1038
$ws.Range("Q33").Value = 15
# row 34 <- original row 58
$ws.Range("D34").Value = 44448
$ws.Range("J34").Value = 85
$ws.Range("K34").Value = 21000
$ws.Range("L34").Value = 22000
$ws.Range("M34").Value = 21529
$ws.Range("N34").Value = '$/malla 15 kilos'
$ws.Range("O34").Value = 'Provincia de Quillota'
$ws.Range("P34").Value = 1435
$ws.Range("Q34").Value = 15
# row 35 <- original row 9
$ws.Range("D35").Value = 44757
$ws.Range("J35").Value = 40
$ws.Range("K35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = 15000
$ws.Range("N35").Value = '$/malla 15 kilos'
$ws.Range("O35").Value = 'Provincia de Quillota'
$ws.Range("P35").Value = 1000
$ws.Range("Q35").Value = 15
# row 36 <- original row 51
$ws.Range("D36").Value = 44824
$ws.Range("J36").Value = 20
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = 15000
$ws.Range("N36").Value = '$/malla 15 kilos'
$ws.Range("O36").Value = 'Provincia de Quillota'
$ws.Range("P36").Value = 1000
$ws.Range("Q36").Value = 15
# row 37 <- original row 84
$ws.Range("D37").Value = 44333
$ws.Range("J37").Value = 35
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 15000
$ws.Range("N37").Value = '$/malla 15 kilos'
$ws.Range("O37").Value = 'Provincia de Quillota'
$ws.Range("P37").Value = 1000
$ws.Range("Q37").Value = 15
# row 38 <- original row 76
$ws.Range("D38").Value = 44754
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = 15000
$ws.Range("N38").Value = '$/malla 15 kilos'
$ws.Range("O38").Value = 'Provincia de Quillota'
$ws.Range("P38").Value = 1000
$ws.Range("Q38").Value = 15
# row 40 <- original row 70
$ws.Range("D40").Value = 44377
$ws.Range("J40").Value = 80
$ws.Range("K40").Value = 18000
$ws.Range("L40").Value = 19000
$ws.Range("M40").Value = 18500
$ws.Range("N40").Value = '$/malla 15 kilos'
$ws.Range("O40").Value = 'Provincia de Quillota'
$ws.Range("P40").Value = 1233
$ws.Range("Q40").Value = 15
# row 41 <- original row 47
$ws.Range("D41").Value = 44314
$ws.Range("J41").Value = 45
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = 15000
$ws.Range("N41").Value = '$/malla 15 kilos'
$ws.Range("O41").Value = 'Provincia de Quillota'
$ws.Range("P41").Value = 1000
$ws.Range("Q41").Value = 15
# row 42 <- original row 19
$ws.Range("D42").Value = 44313
$ws.Range("J42").Value = 40
$ws.Range("K42").Value = 14000
$ws.Range("L42").Value = 14000
$ws.Range("M42").Value = 14000
$ws.Range("N42").Value = '$/malla 15 kilos'
$ws.Range("O42").Value = 'Provincia de Quillota'
$ws.Range("P42").Value = 933
$ws.Range("Q42").Value = 15
# row 43 <- original row 18
$ws.Range("D43").Value = 44326
$ws.Range("J43").Value = 45
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = 15000
$ws.Range("N43").Value = '$/malla 15 kilos'
$ws.Range("O43").Value = 'Provincia de Quillota'
$ws.Range("P43").Value = 1000
$ws.Range("Q43").Value = 15
# row 44 <- original row 14
$ws.Range("D44").Value = 44344
$ws.Range("J44").Value = 40
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = 20000
$ws.Range("N44").Value = '$/malla 15 kilos'
$ws.Range("O44").Value = 'Provincia de Quillota'
$ws.Range("P44").Value = 1333
$ws.Range("Q44").Value = 15
# row 45 <- original row 13
$ws.Range("D45").Value = 44819
$ws.Range("J45").Value = 45
$ws.Range("K45").Value = 16000
$ws.Range("L45").Value = 16000
$ws.Range("M45").Value = 16000
$ws.Range("N45").Value = '$/malla 15 kilos'
$ws.Range("O45").Value = 'Provincia de Quillota'
$ws.Range("P45").Value = 1067
$ws.Range("Q45").Value = 15
# row 46 <- original row 21
$ws.Range("D46").Value = 44798
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 14000
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = 14450
$ws.Range("N46").Value = '$/malla 15 kilos'
$ws.Range("O46").Value = 'Provincia de Quillota'
$ws.Range("P46").Value = 963
$ws.Range("Q46").Value = 15
# row 47 <- original row 43
$ws.Range("D47").Value = 44827
$ws.Range("J47").Value = 45
$ws.Range("K47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = 15000
$ws.Range("N47").Value = '$/malla 15 kilos'
$ws.Range("O47").Value = 'Provincia de Quillota'
$ws.Range("P47").Value = 1000
$ws.Range("Q47").Value = 15
# row 48 <- original row 74
$ws.Range("D48").Value = 44726
$ws.Range("J48").Value = 55
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = 15000
$ws.Range("N48").Value = '$/malla 15 kilos'
$ws.Range("O48").Value = 'Provincia de Quillota'
$ws.Range("P48").Value = 1000
$ws.Range("Q48").Value = 15
# row 49 <- original row 33
$ws.Range("D49").Value = 44792
$ws.Range("J49").Value = 50
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = 15000
$ws.Range("N49").Value = '$/malla 15 kilos'
$ws.Range("O49").Value = 'Provincia de Quillota'
$ws.Range("P49").Value = 1000
$ws.Range("Q49").Value = 15
# row 50 <- original row 22
$ws.Range("D50").Value = 44308
$ws.Range("J50").Value = 40
$ws.Range("K50").Value = 16000
$ws.Range("L50").Value = 16000
$ws.Range("M50").Value = 16000
$ws.Range("N50").Value = '$/malla 15 kilos'
$ws.Range("O50").Value = 'Provincia de Quillota'
$ws.Range("P50").Value = 1067
$ws.Range("Q50").Value = 15
# row 51 <- original row 89
$ws.Range("D51").Value = 44750
$ws.Range("J51").Value = 85
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 16000
$ws.Range("M51").Value = 15471
$ws.Range("N51").Value = '$/malla 15 kilos'
$ws.Range("O51").Value = 'Provincia de Quillota'
$ws.Range("P51").Value = 1031
$ws.Range("Q51").Value = 15
# row 52 <- original row 63
$ws.Range("D52").Value = 44831
$ws.Range("J52").Value = 40
$ws.Range("K52").Value = 16000
$ws.Range("L52").Value = 16000
$ws.Range("M52").Value = 16000
$ws.Range("N52").Value = '$/malla 15 kilos'
$ws.Range("O52").Value = 'Provincia de Quillota'
$ws.Range("P52").Value = 1067
$ws.Range("Q52").Value = 15
# row 53 <- original row 4
$ws.Range("D53").Value = 44817
$ws.Range("J53").Value = 85
$ws.Range("K53").Value = 15000
$ws.Range("L53").Value = 16000
$ws.Range("M53").Value = 15529
$ws.Range("N53").Value = '$/malla 15 kilos'
$ws.Range("O53").Value = 'Provincia de Quillota'
$ws.Range("P53").Value = 1035
$ws.Range("Q53").Value = 15
# row 54 <- original row 32
$ws.Range("D54").Value = 44321
$ws.Range("J54").Value = 38
$ws.Range("K54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("M54").Value = 15000
$ws.Range("N54").Value = '$/malla 15 kilos'
$ws.Range("O54").Value = 'Provincia de Quillota'
$ws.Range("P54").Value = 1000
$ws.Range("Q54").Value = 15
# row 55 <- original row 97
$ws.Range("D55").Value = 44777
$ws.Range("J55").Value = 85
$ws.Range("K55").Value = 14500
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = 14735
$ws.Range("N55").Value = '$/malla 15 kilos'
$ws.Range("O55").Value = 'Provincia de Quillota'
$ws.Range("P55").Value = 982
$ws.Range("Q55").Value = 15
# row 56 <- original row 68
$ws.Range("D56").Value = 44816
$ws.Range("J56").Value = 60
$ws.Range("K56").Value = 16000
$ws.Range("L56").Value = 16000
$ws.Range("M56").Value = 16000
$ws.Range("N56").Value = '$/malla 15 kilos'
$ws.Range("O56").Value = 'Provincia de Quillota'
$ws.Range("P56").Value = 1067
$ws.Range("Q56").Value = 15
# row 57 <- original row 61
$ws.Range("D57").Value = 44797
$ws.Range("J57").Value = 40
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 15000
$ws.Range("M57").Value = 15000
$ws.Range("N57").Value = '$/malla 15 kilos'
$ws.Range("O57").Value = 'Provincia de Quillota'
$ws.Range("P57").Value = 1000
$ws.Range("Q57").Value = 15
# row 58 <- original row 28
$ws.Range("D58").Value = 44312
$ws.Range("J58").Value = 80
$ws.Range("K58").Value = 13000
$ws.Range("L58").Value = 14000
$ws.Range("M58").Value = 13562
$ws.Range("N58").Value = '$/malla 15 kilos'
$ws.Range("O58").Value = 'Provincia de Quillota'
$ws.Range("P58").Value = 904
$ws.Range("Q58").Value = 15
# row 59 <- original row 103
$ws.Range("D59").Value = 44336
$ws.Range("J59").Value = 65
$ws.Range("K59").Value = 14000
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = 14462
$ws.Range("N59").Value = '$/malla 15 kilos'
$ws.Range("O59").Value = 'Provincia de Quillota'
$ws.Range("P59").Value = 964
$ws.Range("Q59").Value = 15
# row 60 <- original row 96
$ws.Range("D60").Value = 44761
$ws.Range("J60").Value = 50
$ws.Range("K60").Value = 15000
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = 15000
$ws.Range("N60").Value = '$/malla 15 kilos'
$ws.Range("O60").Value = 'Provincia de Quillota'
$ws.Range("P60").Value = 1000
$ws.Range("Q60").Value = 15
# row 61 <- original row 26
$ws.Range("D61").Value = 44785
$ws.Range("J61").Value = 85
$ws.Range("K61").Value = 14000
$ws.Range("L61").Value = 15000
$ws.Range("M61").Value = 14471
$ws.Range("N61").Value = '$/malla 15 kilos'
$ws.Range("O61").Value = 'Provincia de Quillota'
$ws.Range("P61").Value = 965
$ws.Range("Q61").Value = 15
# row 62 <- original row 52
$ws.Range("D62").Value = 44764
$ws.Range("J62").Value = 45
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = 15000
$ws.Range("N62").Value = '$/malla 15 kilos'
$ws.Range("O62").Value = 'Provincia de Quillota'
$ws.Range("P62").Value = 1000
$ws.Range("Q62").Value = 15
# row 63 <- original row 2
$ws.Range("D63").Value = 44826
$ws.Range("J63").Value = 50
$ws.Range("K63").Value = 15000
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = 15000
$ws.Range("N63").Value = '$/malla 15 kilos'
$ws.Range("O63").Value = 'Provincia de Quillota'
$ws.Range("P63").Value = 1000
$ws.Range("Q63").Value = 15
# row 64 <- original row 12
$ws.Range("D64").Value = 44762
$ws.Range("J64").Value = 80
$ws.Range("K64").Value = 14000
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = 14500
$ws.Range("N64").Value = '$/malla 15 kilos'
$ws.Range("O64").Value = 'Provincia de Quillota'
$ws.Range("P64").Value = 967
$ws.Range("Q64").Value = 15
# row 65 <- original row 93
$ws.Range("D65").Value = 44334
$ws.Range("J65").Value = 50
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 14000
$ws.Range("M65").Value = 14000
$ws.Range("N65").Value = '$/malla 15 kilos'
$ws.Range("O65").Value = 'Provincia de Quillota'
$ws.Range("P65").Value = 933
$ws.Range("Q65").Value = 15
# row 66 <- original row 57
$ws.Range("D66").Value = 44736
$ws.Range("J66").Value = 82
$ws.Range("K66").Value = 16000
$ws.Range("L66").Value = 17000
$ws.Range("M66").Value = 16488
$ws.Range("N66").Value = '$/malla 15 kilos'
$ws.Range("O66").Value = 'Provincia de Quillota'
$ws.Range("P66").Value = 1099
$ws.Range("Q66").Value = 15
# row 67 <- original row 45
$ws.Range("D67").Value = 44760
$ws.Range("J67").Value = 105
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 16000
$ws.Range("M67").Value = 15524
$ws.Range("N67").Value = '$/malla 15 kilos'
$ws.Range("O67").Value = 'Provincia de Quillota'
$ws.Range("P67").Value = 1035
$ws.Range("Q67").Value = 15
# row 68 <- original row 36
$ws.Range("D68").Value = 44841
$ws.Range("J68").Value = 38
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("M68").Value = 15000
$ws.Range("N68").Value = '$/malla 15 kilos'
$ws.Range("O68").Value = 'Provincia de Quillota'
$ws.Range("P68").Value = 1000
$ws.Range("Q68").Value = 15
# row 69 <- original row 37
$ws.Range("D69").Value = 44809
$ws.Range("J69").Value = 105
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 16000
$ws.Range("M69").Value = 15476
$ws.Range("N69").Value = '$/malla 15 kilos'
$ws.Range("O69").Value = 'Provincia de Quillota'
$ws.Range("P69").Value = 1032
$ws.Range("Q69").Value = 15
# row 70 <- original row 99
$ws.Range("D70").Value = 44839
$ws.Range("J70").Value = 40
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = 15000
$ws.Range("N70").Value = '$/malla 15 kilos'
$ws.Range("O70").Value = 'Provincia de Quillota'
$ws.Range("P70").Value = 1000
$ws.Range("Q70").Value = 15
# row 71 <- original row 49
$ws.Range("D71").Value = 44818
$ws.Range("J71").Value = 58
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 16000
$ws.Range("M71").Value = 16000
$ws.Range("N71").Value = '$/malla 15 kilos'
$ws.Range("O71").Value = 'Provincia de Quillota'
$ws.Range("P71").Value = 1067
$ws.Range("Q71").Value = 15
# row 72 <- original row 88
$ws.Range("D72").Value = 44806
$ws.Range("J72").Value = 45
$ws.Range("K72").Value = 16000
$ws.Range("L72").Value = 16000
$ws.Range("M72").Value = 16000
$ws.Range("N72").Value = '$/malla 15 kilos'
$ws.Range("O72").Value = 'Provincia de Quillota'
$ws.Range("P72").Value = 1067
$ws.Range("Q72").Value = 15
# row 73 <- original row 41
$ws.Range("D73").Value = 44316
$ws.Range("J73").Value = 45
$ws.Range("K73").Value = 14000
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = 14444
$ws.Range("N73").Value = '$/malla 15 kilos'
$ws.Range("O73").Value = 'Provincia de Quillota'
$ws.Range("P73").Value = 963
$ws.Range("Q73").Value = 15
# row 74 <- original row 80
$ws.Range("D74").Value = 44832
$ws.Range("J74").Value = 40
$ws.Range("K74").Value = 16000
$ws.Range("L74").Value = 16000
$ws.Range("M74").Value = 16000
$ws.Range("N74").Value = '$/malla 15 kilos'
$ws.Range("O74").Value = 'Provincia de Quillota'
$ws.Range("P74").Value = 1067
$ws.Range("Q74").Value = 15
# row 75 <- original row 10
$ws.Range("D75").Value = 44789
$ws.Range("J75").Value = 40
$ws.Range("K75").Value = 15000
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = 15000
$ws.Range("N75").Value = '$/malla 15 kilos'
$ws.Range("O75").Value = 'Provincia de Quillota'
$ws.Range("P75").Value = 1000
$ws.Range("Q75").Value = 15
# row 76 <- original row 54
$ws.Range("D76").Value = 44729
$ws.Range("J76").Value = 85
$ws.Range("K76").Value = 16000
$ws.Range("L76").Value = 17000
$ws.Range("M76").Value = 16529
$ws.Range("N76").Value = '$/malla 15 kilos'
$ws.Range("O76").Value = 'Provincia de Quillota'
$ws.Range("P76").Value = 1102
$ws.Range("Q76").Value = 15
# row 77 <- original row 30
$ws.Range("D77").Value = 44776
$ws.Range("J77").Value = 105
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 15500
$ws.Range("M77").Value = 15238
$ws.Range("N77").Value = '$/malla 15 kilos'
$ws.Range("O77").Value = 'Provincia de Quillota'
$ws.Range("P77").Value = 1016
$ws.Range("Q77").Value = 15
# row 78 <- original row 102
$ws.Range("D78").Value = 44747
$ws.Range("J78").Value = 40
$ws.Range("K78").Value = 16000
$ws.Range("L78").Value = 16000
$ws.Range("M78").Value = 16000
$ws.Range("N78").Value = '$/malla 15 kilos'
$ws.Range("O78").Value = 'Provincia de Quillota'
$ws.Range("P78").Value = 1067
$ws.Range("Q78").Value = 15
# row 79 <- original row 8
$ws.Range("D79").Value = 44328
$ws.Range("J79").Value = 38
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = 15000
$ws.Range("N79").Value = '$/malla 15 kilos'
$ws.Range("O79").Value = 'Provincia de Quillota'
$ws.Range("P79").Value = 1000
$ws.Range("Q79").Value = 15
# row 80 <- original row 29
$ws.Range("D80").Value = 44790
$ws.Range("J80").Value = 40
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = 15000
$ws.Range("N80").Value = '$/malla 15 kilos'
$ws.Range("O80").Value = 'Provincia de Quillota'
$ws.Range("P80").Value = 1000
$ws.Range("Q80").Value = 15
# row 81 <- original row 48
$ws.Range("D81").Value = 44775
$ws.Range("J81").Value = 93
$ws.Range("K81").Value = 14000
$ws.Range("L81").Value = 15000
$ws.Range("M81").Value = 14516
$ws.Range("N81").Value = '$/malla 15 kilos'
$ws.Range("O81").Value = 'Provincia de Quillota'
$ws.Range("P81").Value = 968
$ws.Range("Q81").Value = 15
# row 82 <- original row 72
$ws.Range("D82").Value = 44714
$ws.Range("J82").Value = 100
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 15500
$ws.Range("M82").Value = 15250
$ws.Range("N82").Value = '$/malla 15 kilos'
$ws.Range("O82").Value = 'Provincia de Quillota'
$ws.Range("P82").Value = 1017
$ws.Range("Q82").Value = 15
# row 83 <- original row 75
$ws.Range("D83").Value = 44845
$ws.Range("J83").Value = 42
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = 15000
$ws.Range("N83").Value = '$/malla 15 kilos'
$ws.Range("O83").Value = 'Provincia de Quillota'
$ws.Range("P83").Value = 1000
$ws.Range("Q83").Value = 15
# row 84 <- original row 94
$ws.Range("D84").Value = 44763
$ws.Range("J84").Value = 80
$ws.Range("K84").Value = 14000
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = 14500
$ws.Range("N84").Value = '$/malla 15 kilos'
$ws.Range("O84").Value = 'Provincia de Quillota'
$ws.Range("P84").Value = 967
$ws.Range("Q84").Value = 15
# row 85 <- original row 3
$ws.Range("D85").Value = 44320
$ws.Range("J85").Value = 40
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = 15000
$ws.Range("N85").Value = '$/malla 15 kilos'
$ws.Range("O85").Value = 'Provincia de Quillota'
$ws.Range("P85").Value = 1000
$ws.Range("Q85").Value = 15
# row 86 <- original row 35
$ws.Range("D86").Value = 44330
$ws.Range("J86").Value = 30
$ws.Range("K86").Value = 15000
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = 15000
$ws.Range("N86").Value = '$/malla 15 kilos'
$ws.Range("O86").Value = 'Provincia de Quillota'
$ws.Range("P86").Value = 1000
$ws.Range("Q86").Value = 15
# row 87 <- original row 86
$ws.Range("D87").Value = 44742
$ws.Range("J87").Value = 85
$ws.Range("K87").Value = 15000
$ws.Range("L87").Value = 16000
$ws.Range("M87").Value = 15529
$ws.Range("N87").Value = '$/malla 15 kilos'
$ws.Range("O87").Value = 'Provincia de Quillota'
$ws.Range("P87").Value = 1035
$ws.Range("Q87").Value = 15
# row 88 <- original row 59
$ws.Range("D88").Value = 44715
$ws.Range("J88").Value = 85
$ws.Range("K88").Value = 15000
$ws.Range("L88").Value = 15500
$ws.Range("M88").Value = 15235
$ws.Range("N88").Value = '$/malla 15 kilos'
$ws.Range("O88").Value = 'Provincia de Quillota'
$ws.Range("P88").Value = 1016
$ws.Range("Q88").Value = 15
# row 89 <- original row 17
$ws.Range("D89").Value = 44812
$ws.Range("J89").Value = 45
$ws.Range("K89").Value = 16000
$ws.Range("L89").Value = 16000
$ws.Range("M89").Value = 16000
$ws.Range("N89").Value = '$/malla 15 kilos'
$ws.Range("O89").Value = 'Provincia de Quillota'
$ws.Range("P89").Value = 1067
$ws.Range("Q89").Value = 15
# row 90 <- original row 60
$ws.Range("D90").Value = 44397
$ws.Range("J90").Value = 73
$ws.Range("K90").Value = 21000
$ws.Range("L90").Value = 22000
$ws.Range("M90").Value = 21521
$ws.Range("N90").Value = '$/malla 15 kilos'
$ws.Range("O90").Value = 'Provincia de Quillota'
$ws.Range("P90").Value = 1435
$ws.Range("Q90").Value = 15
# row 91 <- original row 56
$ws.Range("D91").Value = 44767
$ws.Range("J91").Value = 45
$ws.Range("K91").Value = 15000
$ws.Range("L91").Value = 15000
$ws.Range("M91").Value = 15000
$ws.Range("N91").Value = '$/malla 15 kilos'
$ws.Range("O91").Value = 'Provincia de Quillota'
$ws.Range("P91").Value = 1000
$ws.Range("Q91").Value = 15
# row 92 <- original row 42
$ws.Range("D92").Value = 44825
$ws.Range("J92").Value = 85
$ws.Range("K92").Value = 15000
$ws.Range("L92").Value = 15500
$ws.Range("M92").Value = 15265
$ws.Range("N92").Value = '$/malla 15 kilos'
$ws.Range("O92").Value = 'Provincia de Quillota'
$ws.Range("P92").Value = 1018
$ws.Range("Q92").Value = 15
# row 93 <- original row 87
$ws.Range("D93").Value = 44803
$ws.Range("J93").Value = 85
$ws.Range("K93").Value = 15000
$ws.Range("L93").Value = 15500
$ws.Range("M93").Value = 15265
$ws.Range("N93").Value = '$/malla 15 kilos'
$ws.Range("O93").Value = 'Provincia de Quillota'
$ws.Range("P93").Value = 1018
$ws.Range("Q93").Value = 15
# row 94 <- original row 40
$ws.Range("D94").Value = 44370
$ws.Range("J94").Value = 50
$ws.Range("K94").Value = 18000
$ws.Range("L94").Value = 18000
$ws.Range("M94").Value = 18000
$ws.Range("N94").Value = '$/malla 15 kilos'
$ws.Range("O94").Value = 'Provincia de Quillota'
$ws.Range("P94").Value = 1200
$ws.Range("Q94").Value = 15
# row 95 <- original row 16
$ws.Range("D95").Value = 44799
$ws.Range("J95").Value = 55
$ws.Range("K95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("M95").Value = 15000
$ws.Range("N95").Value = '$/malla 15 kilos'
$ws.Range("O95").Value = 'Provincia de Quillota'
$ws.Range("P95").Value = 1000
$ws.Range("Q95").Value = 15
# row 96 <- original row 100
$ws.Range("D96").Value = 44755
$ws.Range("J96").Value = 100
$ws.Range("K96").Value = 15000
$ws.Range("L96").Value = 16000
$ws.Range("M96").Value = 15550
$ws.Range("N96").Value = '$/malla 15 kilos'
$ws.Range("O96").Value = 'Provincia de Quillota'
$ws.Range("P96").Value = 1037
$ws.Range("Q96").Value = 15
# row 97 <- original row 55
$ws.Range("D97").Value = 44323
$ws.Range("J97").Value = 40
$ws.Range("K97").Value = 15000
$ws.Range("L97").Value = 15000
$ws.Range("M97").Value = 15000
$ws.Range("N97").Value = '$/malla 15 kilos'
$ws.Range("O97").Value = 'Provincia de Quillota'
$ws.Range("P97").Value = 1000
$ws.Range("Q97").Value = 15
# row 98 <- original row 31
$ws.Range("D98").Value = 44327
$ws.Range("J98").Value = 35
$ws.Range("K98").Value = 15000
$ws.Range("L98").Value = 15000
$ws.Range("M98").Value = 15000
$ws.Range("N98").Value = '$/malla 15 kilos'
$ws.Range("O98").Value = 'Provincia de Quillota'
$ws.Range("P98").Value = 1000
$ws.Range("Q98").Value = 15
# row 99 <- original row 78
$ws.Range("D99").Value = 44341
$ws.Range("J99").Value = 40
$ws.Range("K99").Value = 15000
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = 15000
$ws.Range("N99").Value = '$/malla 15 kilos'
$ws.Range("O99").Value = 'Provincia de Quillota'
$ws.Range("P99").Value = 1000
$ws.Range("Q99").Value = 15
# row 100 <- original row 27
$ws.Range("D100").Value = 44756
$ws.Range("J100").Value = 50
$ws.Range("K100").Value = 15000
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = 15000
$ws.Range("N100").Value = '$/malla 15 kilos'
$ws.Range("O100").Value = 'Provincia de Quillota'
$ws.Range("P100").Value = 1000
$ws.Range("Q100").Value = 15
# row 101 <- original row 98
$ws.Range("D101").Value = 44810
$ws.Range("J101").Value = 85
$ws.Range("K101").Value = 16000
$ws.Range("L101").Value = 16500
$ws.Range("M101").Value = 16235
$ws.Range("N101").Value = '$/malla 15 kilos'
$ws.Range("O101").Value = 'Provincia de Quillota'
$ws.Range("P101").Value = 1082
$ws.Range("Q101").Value = 15
# row 102 <- original row 5
$ws.Range("D102").Value = 44769
$ws.Range("J102").Value = 85
$ws.Range("K102").Value = 14000
$ws.Range("L102").Value = 15000
$ws.Range("M102").Value = 14471
$ws.Range("N102").Value = '$/malla 15 kilos'
$ws.Range("O102").Value = 'Provincia de Quillota'
$ws.Range("P102").Value = 965
$ws.Range("Q102").Value = 15
# row 103 <- original row 85
$ws.Range("D103").Value = 44725
$ws.Range("J103").Value = 85
$ws.Range("K103").Value = 14000
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 14471
$ws.Range("N103").Value = '$/malla 15 kilos'
$ws.Range("O103").Value = 'Provincia de Quillota'
$ws.Range("P103").Value = 965
$ws.Range("Q103").Value = 15
# row 104 <- original row 15
$ws.Range("D104").Value = 44837
$ws.Range("J104").Value = 40
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("M104").Value = 15000
$ws.Range("N104").Value = '$/malla 15 kilos'
$ws.Range("O104").Value = 'Provincia de Quillota'
$ws.Range("P104").Value = 1000
$ws.Range("Q104").Value = 15
